$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Row 2 / column B ("ParticipantsTab" query) is rewritten with a corrected
# Cypher query that uses OPTIONAL MATCH so participants without matching
# genomic_info/file rows are not silently dropped by the primary-diagnosis
# filter ("cds Primary diagnosis fixed").
$newQuery = @"
MATCH (p:participant)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
OPTIONAL MATCH (p)<--(diag:diagnosis)
OPTIONAL MATCH (samp)<--(f:file)
OPTIONAL MATCH (f)<--(g:genomic_info)
WITH s, p, samp, f, g, diag
WHERE diag.primary_diagnosis in ['Acute monoblastic leukemia']
WITH p
OPTIONAL MATCH (p)-->(s:study)
OPTIONAL MATCH (samp:sample)-->(p)
WITH s, p, apoc.coll.sort(collect(distinct samp.sample_id)) as samp
RETURN 
coalesce(p.participant_id,'') as ``Participant ID``,
coalesce(s.study_name, '') as ``Study Name``,
coalesce(s.phs_accession,'') as ``Accession``,
coalesce(p.gender,'') as ``Gender``,
coalesce(apoc.text.join(samp, ','), '') as ``Samples``
ORDER BY p.participant_id
LIMIT 100
"@

$ws.Range("B2").Value = $newQuery

# The longer replacement text wraps onto more lines, so the row grows taller.
$ws.Rows.Item(2).RowHeight = 330.75

# Selection moved to D3 in the saved workbook.
$ws.Range("D3").Select()
